$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Havel Partners"
$ws.Range("B2").Value = "17sec"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "1"
$ws.Range("C2").NumberFormat = "General"

$ws.Range("A3").Value = "Samvad Partners"
$ws.Range("B3").Value = "13sec"

$ws.Range("A4").Value = "Byrne Wallace"
$ws.Range("B4").Value = "10sec"

$ws.Range("A5").Value = "Anand And Anand"
$ws.Range("B5").Value = "13sec"

$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
